$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# Update the "last updated" timestamp in A1
$ws.Range("A1").Value = "Datos actualizados a 19 de Septiembre de 2020 a las 15:44"

# Update country rows whose rank / stats changed with the new data refresh.
# (rows are re-sorted descending by total cases; a country that overtakes its
# neighbour gets fresh numbers while the other keeps its old numbers one row down)
$ws.Cells.Item(4, 1).Value = "Estados Unidos"
$ws.Cells.Item(4, 2).Value = 6928412
$ws.Cells.Item(4, 3).Value = 2471
$ws.Cells.Item(4, 4).Value = 4192774
$ws.Cells.Item(4, 5).Value = 2532450
$ws.Cells.Item(4, 6).Value = 0
$ws.Cells.Item(4, 7).Value = 17
$ws.Cells.Item(4, 8).Value = 203188

$ws.Cells.Item(5, 1).Value = "India"
$ws.Cells.Item(5, 2).Value = 5322415
$ws.Cells.Item(5, 3).Value = 16940
$ws.Cells.Item(5, 4).Value = 4220511
$ws.Cells.Item(5, 5).Value = 1016194
$ws.Cells.Item(5, 6).Value = 0
$ws.Cells.Item(5, 7).Value = 85
$ws.Cells.Item(5, 8).Value = 85710

$ws.Cells.Item(19, 1).Value = "Arabia Saudita"
$ws.Cells.Item(19, 2).Value = 329271
$ws.Cells.Item(19, 3).Value = 551
$ws.Cells.Item(19, 4).Value = 309430
$ws.Cells.Item(19, 5).Value = 15383
$ws.Cells.Item(19, 6).Value = 0
$ws.Cells.Item(19, 7).Value = 28
$ws.Cells.Item(19, 8).Value = 4458

$ws.Cells.Item(25, 1).Value = "Alemania"
$ws.Cells.Item(25, 2).Value = 271480
$ws.Cells.Item(25, 3).Value = 236
$ws.Cells.Item(25, 4).Value = 243000
$ws.Cells.Item(25, 5).Value = 19014
$ws.Cells.Item(25, 6).Value = 0
$ws.Cells.Item(25, 7).Value = 2
$ws.Cells.Item(25, 8).Value = 9466

$ws.Cells.Item(39, 1).Value = "Kuwait"
$ws.Cells.Item(39, 2).Value = 99049
$ws.Cells.Item(39, 3).Value = 521
$ws.Cells.Item(39, 4).Value = 89498
$ws.Cells.Item(39, 5).Value = 8970
$ws.Cells.Item(39, 6).Value = 0
$ws.Cells.Item(39, 7).Value = 1
$ws.Cells.Item(39, 8).Value = 581

$ws.Cells.Item(41, 1).Value = "Paises Bajos"
$ws.Cells.Item(41, 2).Value = 91934
$ws.Cells.Item(41, 3).Value = 1887
$ws.Cells.Item(41, 4).Value = 0
$ws.Cells.Item(41, 5).Value = 0
$ws.Cells.Item(41, 6).Value = 0
$ws.Cells.Item(41, 7).Value = 2
$ws.Cells.Item(41, 8).Value = 6275

$ws.Cells.Item(42, 1).Value = "Oman"
$ws.Cells.Item(42, 2).Value = 91753
$ws.Cells.Item(42, 3).Value = 0
$ws.Cells.Item(42, 4).Value = 84648
$ws.Cells.Item(42, 5).Value = 6287
$ws.Cells.Item(42, 6).Value = 0
$ws.Cells.Item(42, 7).Value = 0
$ws.Cells.Item(42, 8).Value = 818

$ws.Cells.Item(51, 1).Value = "Portugal"
$ws.Cells.Item(51, 2).Value = 68025
$ws.Cells.Item(51, 3).Value = 849
$ws.Cells.Item(51, 4).Value = 45404
$ws.Cells.Item(51, 5).Value = 20722
$ws.Cells.Item(51, 6).Value = 0
$ws.Cells.Item(51, 7).Value = 5
$ws.Cells.Item(51, 8).Value = 1899

$ws.Cells.Item(52, 1).Value = "Etiopia"
$ws.Cells.Item(52, 2).Value = 67515
$ws.Cells.Item(52, 3).Value = 0
$ws.Cells.Item(52, 4).Value = 27638
$ws.Cells.Item(52, 5).Value = 38805
$ws.Cells.Item(52, 6).Value = 0
$ws.Cells.Item(52, 7).Value = 0
$ws.Cells.Item(52, 8).Value = 1072

$ws.Cells.Item(54, 1).Value = "Barein"
$ws.Cells.Item(54, 2).Value = 63879
$ws.Cells.Item(54, 3).Value = 0
$ws.Cells.Item(54, 4).Value = 56700
$ws.Cells.Item(54, 5).Value = 6958
$ws.Cells.Item(54, 6).Value = 0
$ws.Cells.Item(54, 7).Value = 1
$ws.Cells.Item(54, 8).Value = 221

$ws.Cells.Item(67, 1).Value = "Azerbaiyan"
$ws.Cells.Item(67, 2).Value = 39042
$ws.Cells.Item(67, 3).Value = 148
$ws.Cells.Item(67, 4).Value = 36601
$ws.Cells.Item(67, 5).Value = 1867
$ws.Cells.Item(67, 6).Value = 0
$ws.Cells.Item(67, 7).Value = 2
$ws.Cells.Item(67, 8).Value = 574

$ws.Cells.Item(68, 1).Value = "Afganistan"
$ws.Cells.Item(68, 2).Value = 38919
$ws.Cells.Item(68, 3).Value = 36
$ws.Cells.Item(68, 4).Value = 32576
$ws.Cells.Item(68, 5).Value = 4906
$ws.Cells.Item(68, 6).Value = 0
$ws.Cells.Item(68, 7).Value = 0
$ws.Cells.Item(68, 8).Value = 1437

$ws.Cells.Item(70, 1).Value = "Kenia"
$ws.Cells.Item(70, 2).Value = 36829
$ws.Cells.Item(70, 3).Value = 105
$ws.Cells.Item(70, 4).Value = 23777
$ws.Cells.Item(70, 5).Value = 12406
$ws.Cells.Item(70, 6).Value = 0
$ws.Cells.Item(70, 7).Value = 0
$ws.Cells.Item(70, 8).Value = 646

$ws.Cells.Item(71, 1).Value = "Estado de Palestina"
$ws.Cells.Item(71, 2).Value = 35003
$ws.Cells.Item(71, 3).Value = 602
$ws.Cells.Item(71, 4).Value = 23446
$ws.Cells.Item(71, 5).Value = 11304
$ws.Cells.Item(71, 6).Value = 0
$ws.Cells.Item(71, 7).Value = 3
$ws.Cells.Item(71, 8).Value = 253

$ws.Cells.Item(72, 1).Value = "Serbia"
$ws.Cells.Item(72, 2).Value = 32840
$ws.Cells.Item(72, 3).Value = 83
$ws.Cells.Item(72, 4).Value = 31411
$ws.Cells.Item(72, 5).Value = 689
$ws.Cells.Item(72, 6).Value = 0
$ws.Cells.Item(72, 7).Value = 1
$ws.Cells.Item(72, 8).Value = 740

$ws.Cells.Item(77, 1).Value = "Libia"
$ws.Cells.Item(77, 2).Value = 27234
$ws.Cells.Item(77, 3).Value = 796
$ws.Cells.Item(77, 4).Value = 14679
$ws.Cells.Item(77, 5).Value = 12119
$ws.Cells.Item(77, 6).Value = 0
$ws.Cells.Item(77, 7).Value = 18
$ws.Cells.Item(77, 8).Value = 436

$ws.Cells.Item(78, 1).Value = "Australia"
$ws.Cells.Item(78, 2).Value = 26885
$ws.Cells.Item(78, 3).Value = 24
$ws.Cells.Item(78, 4).Value = 23962
$ws.Cells.Item(78, 5).Value = 2079
$ws.Cells.Item(78, 6).Value = 0
$ws.Cells.Item(78, 7).Value = 7
$ws.Cells.Item(78, 8).Value = 844

$ws.Cells.Item(79, 1).Value = "Bosnia y Herzegovina"
$ws.Cells.Item(79, 2).Value = 25217
$ws.Cells.Item(79, 3).Value = 320
$ws.Cells.Item(79, 4).Value = 17489
$ws.Cells.Item(79, 5).Value = 6970
$ws.Cells.Item(79, 6).Value = 0
$ws.Cells.Item(79, 7).Value = 6
$ws.Cells.Item(79, 8).Value = 758

$ws.Cells.Item(81, 1).Value = "Dinamarca"
$ws.Cells.Item(81, 2).Value = 22436
$ws.Cells.Item(81, 3).Value = 589
$ws.Cells.Item(81, 4).Value = 17316
$ws.Cells.Item(81, 5).Value = 4485
$ws.Cells.Item(81, 6).Value = 0
$ws.Cells.Item(81, 7).Value = 0
$ws.Cells.Item(81, 8).Value = 635

$ws.Cells.Item(101, 1).Value = "Tayikistan"
$ws.Cells.Item(101, 2).Value = 9303
$ws.Cells.Item(101, 3).Value = 44
$ws.Cells.Item(101, 4).Value = 8066
$ws.Cells.Item(101, 5).Value = 1164
$ws.Cells.Item(101, 6).Value = 0
$ws.Cells.Item(101, 7).Value = 0
$ws.Cells.Item(101, 8).Value = 73

$ws.Cells.Item(140, 1).Value = "Sri Lanka"
$ws.Cells.Item(140, 2).Value = 3282
$ws.Cells.Item(140, 3).Value = 1
$ws.Cells.Item(140, 4).Value = 3070
$ws.Cells.Item(140, 5).Value = 199
$ws.Cells.Item(140, 6).Value = 0
$ws.Cells.Item(140, 7).Value = 0
$ws.Cells.Item(140, 8).Value = 13

$ws.Cells.Item(164, 1).Value = "Polinesia Francesa"
$ws.Cells.Item(164, 2).Value = 1209
$ws.Cells.Item(164, 3).Value = 98
$ws.Cells.Item(164, 4).Value = 966
$ws.Cells.Item(164, 5).Value = 241
$ws.Cells.Item(164, 6).Value = 0
$ws.Cells.Item(164, 7).Value = 0
$ws.Cells.Item(164, 8).Value = 2

$ws.Cells.Item(165, 1).Value = "Niger"
$ws.Cells.Item(165, 2).Value = 1183
$ws.Cells.Item(165, 3).Value = 0
$ws.Cells.Item(165, 4).Value = 1104
$ws.Cells.Item(165, 5).Value = 10
$ws.Cells.Item(165, 6).Value = 0
$ws.Cells.Item(165, 7).Value = 0
$ws.Cells.Item(165, 8).Value = 69

$ws.Cells.Item(166, 1).Value = "Republica del Chad"
$ws.Cells.Item(166, 2).Value = 1147
$ws.Cells.Item(166, 3).Value = 0
$ws.Cells.Item(166, 4).Value = 966
$ws.Cells.Item(166, 5).Value = 100
$ws.Cells.Item(166, 6).Value = 0
$ws.Cells.Item(166, 7).Value = 0
$ws.Cells.Item(166, 8).Value = 81

$ws.Cells.Item(167, 1).Value = "Martinica"
$ws.Cells.Item(167, 2).Value = 1122
$ws.Cells.Item(167, 3).Value = 0
$ws.Cells.Item(167, 4).Value = 98
$ws.Cells.Item(167, 5).Value = 1006
$ws.Cells.Item(167, 6).Value = 0
$ws.Cells.Item(167, 7).Value = 0
$ws.Cells.Item(167, 8).Value = 18

$ws.Cells.Item(179, 1).Value = "Islas Feroe"
$ws.Cells.Item(179, 2).Value = 431
$ws.Cells.Item(179, 3).Value = 1
$ws.Cells.Item(179, 4).Value = 412
$ws.Cells.Item(179, 5).Value = 19
$ws.Cells.Item(179, 6).Value = 0
$ws.Cells.Item(179, 7).Value = 0
$ws.Cells.Item(179, 8).Value = 0

$ws.Cells.Item(193, 1).Value = "Brunei"
$ws.Cells.Item(193, 2).Value = 145
$ws.Cells.Item(193, 3).Value = 0
$ws.Cells.Item(193, 4).Value = 142
$ws.Cells.Item(193, 5).Value = 0
$ws.Cells.Item(193, 6).Value = 0
$ws.Cells.Item(193, 7).Value = 0
$ws.Cells.Item(193, 8).Value = 3

$ws.Cells.Item(195, 1).Value = "Liechtenstein"
$ws.Cells.Item(195, 2).Value = 113
$ws.Cells.Item(195, 3).Value = 1
$ws.Cells.Item(195, 4).Value = 109
$ws.Cells.Item(195, 5).Value = 3
$ws.Cells.Item(195, 6).Value = 0
$ws.Cells.Item(195, 7).Value = 0
$ws.Cells.Item(195, 8).Value = 1

$ws.Cells.Item(204, 1).Value = "Santa Lucia"
$ws.Cells.Item(204, 2).Value = 27
$ws.Cells.Item(204, 3).Value = 0
$ws.Cells.Item(204, 4).Value = 26
$ws.Cells.Item(204, 5).Value = 1
$ws.Cells.Item(204, 6).Value = 0
$ws.Cells.Item(204, 7).Value = 0
$ws.Cells.Item(204, 8).Value = 0

$ws.Cells.Item(205, 1).Value = "Timor Oriental"
$ws.Cells.Item(205, 2).Value = 27
$ws.Cells.Item(205, 3).Value = 0
$ws.Cells.Item(205, 4).Value = 26
$ws.Cells.Item(205, 5).Value = 1
$ws.Cells.Item(205, 6).Value = 0
$ws.Cells.Item(205, 7).Value = 0
$ws.Cells.Item(205, 8).Value = 0

$ws.Cells.Item(214, 1).Value = "Montserrat"
$ws.Cells.Item(214, 2).Value = 13
$ws.Cells.Item(214, 3).Value = 0
$ws.Cells.Item(214, 4).Value = 12
$ws.Cells.Item(214, 5).Value = 0
$ws.Cells.Item(214, 6).Value = 0
$ws.Cells.Item(214, 7).Value = 0
$ws.Cells.Item(214, 8).Value = 1

$ws.Cells.Item(215, 1).Value = "Islas Malvinas"
$ws.Cells.Item(215, 2).Value = 13
$ws.Cells.Item(215, 3).Value = 0
$ws.Cells.Item(215, 4).Value = 13
$ws.Cells.Item(215, 5).Value = 0
$ws.Cells.Item(215, 6).Value = 0
$ws.Cells.Item(215, 7).Value = 0
$ws.Cells.Item(215, 8).Value = 0

